$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 293.10526
$ws.Range("I11").Value = 293.10526
$ws.Range("K11").Value = 293.10526
$ws.Range("M11").Value = -153.10526

$ws.Range("H64").Value = 3858.5
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 3944.6667
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 3944.6667
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -4440.6667

$ws.Range("H67").Value = 3858.5
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 3944.6667
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 3944.6667
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -5660.6667

$ws.Range("H74").Value = 3495.6667
$ws.Range("I74").Value = 3495.6667
$ws.Range("K74").Value = 3495.6667
$ws.Range("M74").Value = -2559.6667

$ws.Range("H77").Value = 3495.6667
$ws.Range("I77").Value = 3495.6667
$ws.Range("K77").Value = 17478.3335
$ws.Range("M77").Value = -12798.3335

$ws.Range("H106").Value = 5353.3335
$ws.Range("I106").Value = 5353.3335
$ws.Range("K106").Value = 5353.3335
$ws.Range("M106").Value = -4722.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3491.3555
$ws.Range("I32").Value = 1359.8334
$ws.Range("J32").Value = 33332.668
$ws.Range("K32").Value = 1359.8334
$ws.Range("L32").Value = 33332.668
$ws.Range("M32").Value = -1072.8334
$ws.Range("N32").Value = -33906.668

$ws.Range("H61").Value = 2762.182
$ws.Range("I61").Value = 2582.625
$ws.Range("K61").Value = 2582.625
$ws.Range("M61").Value = -2370.625

$ws.Range("H122").Value = 2587.8235
$ws.Range("I122").Value = 2557.1428
$ws.Range("J122").Value = 2731
$ws.Range("K122").Value = 7671.428400000001
$ws.Range("L122").Value = 8193
$ws.Range("M122").Value = -5221.428400000001
$ws.Range("N122").Value = -13093

$ws.Range("H132").Value = 6236.4
$ws.Range("I132").Value = 6236
$ws.Range("J132").Value = 6238
$ws.Range("K132").Value = 18708
$ws.Range("L132").Value = 18714
$ws.Range("M132").Value = -16178
$ws.Range("N132").Value = -23774

$ws.Range("H136").Value = 2762.182
$ws.Range("I136").Value = 2582.625
$ws.Range("K136").Value = 7747.875
$ws.Range("M136").Value = -5197.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 20235.5
$ws.Range("J76").Value = 20235.5
$ws.Range("L76").Value = 20235.5
$ws.Range("N76").Value = -20865.5

$ws.Range("H79").Value = 20235.5
$ws.Range("J79").Value = 20235.5
$ws.Range("L79").Value = 20235.5
$ws.Range("N79").Value = -22419.5

$ws.Range("H86").Value = 2128.0625
$ws.Range("I86").Value = 2303.5
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 2303.5
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = -1180.5
$ws.Range("N86").Value = -3146

$ws.Range("H89").Value = 2128.0625
$ws.Range("I89").Value = 2303.5
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 11517.5
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -5901.5
$ws.Range("N89").Value = -15732

$ws.Range("H134").Value = 5379.1333
$ws.Range("J134").Value = 3665.1875
$ws.Range("L134").Value = 10995.5625
$ws.Range("N134").Value = -16065.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1628.8
$ws.Range("I16").Value = 1628.8
$ws.Range("K16").Value = 1628.8
$ws.Range("M16").Value = -1341.8

$ws.Range("H62").Value = 1934
$ws.Range("I62").Value = 1951
$ws.Range("K62").Value = 1951
$ws.Range("M62").Value = -1327

$ws.Range("H65").Value = 1934
$ws.Range("I65").Value = 1951
$ws.Range("K65").Value = 9755
$ws.Range("M65").Value = -6635

$ws.Range("H109").Value = 39000
$ws.Range("J109").Value = 39000
$ws.Range("L109").Value = 39000
$ws.Range("N109").Value = -41080

$ws.Range("H113").Value = 1628.8
$ws.Range("I113").Value = 1628.8
$ws.Range("K113").Value = 1628.8
$ws.Range("M113").Value = 541.2

$ws.Range("H132").Value = 4587.5
$ws.Range("I132").Value = 3402.8
$ws.Range("J132").Value = 4935.9414
$ws.Range("K132").Value = 10208.4
$ws.Range("L132").Value = 14807.8242
$ws.Range("M132").Value = -7678.400000000001
$ws.Range("N132").Value = -19867.8242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1384.4
$ws.Range("I25").Value = 833
$ws.Range("J25").Value = 3590
$ws.Range("K25").Value = 2499
$ws.Range("L25").Value = 10770
$ws.Range("M25").Value = -2330
$ws.Range("N25").Value = -11108

$ws.Range("H30").Value = 1384.4
$ws.Range("I30").Value = 833
$ws.Range("J30").Value = 3590
$ws.Range("K30").Value = 2499
$ws.Range("L30").Value = 10770
$ws.Range("M30").Value = -2397
$ws.Range("N30").Value = -10974

$ws.Range("H108").Value = 602.5
$ws.Range("I108").Value = 602.5
$ws.Range("K108").Value = 1807.5
$ws.Range("M108").Value = 1072.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2288.25
$ws.Range("I22").Value = 1001
$ws.Range("J22").Value = 2717.3333
$ws.Range("K22").Value = 1001
$ws.Range("L22").Value = 2717.3333
$ws.Range("M22").Value = -706
$ws.Range("N22").Value = -3307.3333

$ws.Range("H27").Value = 2288.25
$ws.Range("I27").Value = 1001
$ws.Range("J27").Value = 2717.3333
$ws.Range("K27").Value = 1001
$ws.Range("L27").Value = 2717.3333
$ws.Range("M27").Value = -894
$ws.Range("N27").Value = -2931.3333

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H136").Value = 3955.5
$ws.Range("I136").Value = 3282
$ws.Range("J136").Value = 4629
$ws.Range("K136").Value = 9846
$ws.Range("L136").Value = 13887
$ws.Range("M136").Value = -7296
$ws.Range("N136").Value = -18987

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H126").Value = 39599.645
$ws.Range("I126").Value = 35366.582
$ws.Range("K126").Value = 106099.746
$ws.Range("M126").Value = -103629.746

$ws.Range("H132").Value = 2321.3635
$ws.Range("I132").Value = 798.6
$ws.Range("J132").Value = 3590.3333
$ws.Range("K132").Value = 2395.8
$ws.Range("L132").Value = 10770.9999
$ws.Range("M132").Value = 134.1999999999998
$ws.Range("N132").Value = -15830.9999

$ws.Range("H136").Value = 5274.421
$ws.Range("I136").Value = 3469.7
$ws.Range("K136").Value = 10409.1
$ws.Range("M136").Value = -7859.7859999999999
